$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M7").Value = 1299.08
$ws1.Range("H15").Value = 963
$ws1.Range("M15").Value = 230.94
$ws1.Range("M22").Value = "10 de 20"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F7").Value = 2898.45
$ws2.Range("F15").Value = 1319.85
$ws2.Range("F22").Value = 51008.12

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Columns.Item(5).ColumnWidth = 22.166666666666668

$ws3.Range("D7").Value = 3339
$ws3.Range("E7").Value = -939
$ws3.Range("F7").Value = 1.39125

$ws3.Range("D16").Value = 41174.07
$ws3.Range("E16").Value = 3092.169999999998
$ws3.Range("F16").Value = 0.9301460887574821

$ws3.Range("D19").Value = 51008.12
$ws3.Range("E19").Value = 14369.87762291768
$ws3.Range("F19").Value = 0.7802031548014182
